$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the stale hyperlinks so we can rebuild them (with explicit
#    display text) further below.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2) Update the data rows 2-5 (existing contacts) with the new values.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "ahostess-test"
$ws.Range("B2").Value = "abc"
$ws.Range("C2").Value = "ahostess-test@test.com"
$ws.Range("D2").Value = "973 BRAHMS CT"
$ws.Range("E2").Value = "TROY"
$ws.Range("F2").Value = "Michigan"
$ws.Range("G2").Value = 48085
$ws.Range("H2").Value = "blackdress19"

$ws.Range("A3").Value = "bcohost-test"
$ws.Range("B3").Value = "abc"
$ws.Range("C3").Value = "bcohost-test@test.com"
$ws.Range("D3").Value = "974 BRAHMS CT"
$ws.Range("E3").Value = "TROY"
$ws.Range("F3").Value = "Michigan"
$ws.Range("G3").Value = 48085
$ws.Range("H3").Value = "blackdress19"

$ws.Range("A4").Value = "guest1-test"
$ws.Range("B4").Value = "abc"
$ws.Range("C4").Value = "guest1-test@test.com"
$ws.Range("D4").Value = "975 BRAHMS CT"
$ws.Range("E4").Value = "TROY"
$ws.Range("F4").Value = "Michigan"
$ws.Range("G4").Value = 48085
$ws.Range("H4").Value = "blackdress19"

$ws.Range("A5").Value = "guest2-test"
$ws.Range("B5").Value = "abc"
$ws.Range("C5").Value = "guest2-test@test.com"
$ws.Range("D5").Value = "976 BRAHMS CT"
$ws.Range("E5").Value = "TROY"
$ws.Range("F5").Value = "Michigan"
$ws.Range("G5").Value = 48085
$ws.Range("H5").Value = "blackdress19"

# ---------------------------------------------------------------------
# 3) Add the brand-new 6th contact row, copying formats from row 5 so
#    the new row inherits the same column styles (s=4 / s=2 / s=3).
# ---------------------------------------------------------------------
$ws.Range("A5:H5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows.Item(6).RowHeight = 15

$ws.Range("A6").Value = "guest3-test"
$ws.Range("B6").Value = "abc"
$ws.Range("C6").Value = "guest3-test@test.com"
$ws.Range("D6").Value = "976 BRAHMS CT"
$ws.Range("E6").Value = "TROY"
$ws.Range("F6").Value = "Michigan"
$ws.Range("G6").Value = 48085
$ws.Range("H6").Value = "blackdress19"

# ---------------------------------------------------------------------
# 4) Rebuild the hyperlinks (rId1-4 keep pointing at the legacy "test3"
#    mailboxes that were already wired up before; rId5 is the brand
#    new 6th contact's hyperlink).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ahostess-test3@test.com", [Type]::Missing, [Type]::Missing, "ahostess-test3@test.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:bcohost-test3@test.com", [Type]::Missing, [Type]::Missing, "bcohost-test3@test.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:guest1-test3@test.com", [Type]::Missing, [Type]::Missing, "guest1-test3@test.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:guest2-test3@test.com", [Type]::Missing, [Type]::Missing, "guest2-test3@test.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:guest3-test1@test.com", [Type]::Missing, [Type]::Missing, "guest3-test1@test.com")

# Hyperlinks.Add overwrites the cell text (TextToDisplay) and the cell
# style when it attaches to a range, so put the real contact e-mail
# text and the Hyperlink cell style back now that the links exist.
$ws.Range("C2").Value = "ahostess-test@test.com"
$ws.Range("C3").Value = "bcohost-test@test.com"
$ws.Range("C4").Value = "guest1-test@test.com"
$ws.Range("C5").Value = "guest2-test@test.com"
$ws.Range("C6").Value = "guest3-test@test.com"

$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"
$ws.Range("C6").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 5) Match the final cursor position recorded in the worksheet view.
# ---------------------------------------------------------------------
[void]$ws.Range("F15").Select()
